# "add nptg add overwrite option"
# The NaPTAN request template's "Stops" sheet had a redundant "Street"
# column (with its data-validation list restricted to the placeholder "-")
# that duplicated the NptgLocality-driven workflow. Remove that column so
# the sheet goes from AtcoCode..AdministrativeAreaRef across columns A:R
# down to A:Q, shifting StopType..AdministrativeAreaRef one column left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stops")

# Delete the whole "Street" column (E) - removes the header cell, the
# sample data cell, and shifts every column after it (data validations
# included) one place to the left.
$ws.Columns("E:E").Delete() | Out-Null

# Leave the cursor where the author left it when they saved the file.
$ws.Range("F9").Select() | Out-Null
